# "fixed dv + priority"
#
# The sheet is a drone delivery queue: col A is a countdown/priority number,
# B/C/D ("From"/"To"/"Drone") are stored as text (shared strings) even
# though they look numeric. The fix re-numbers the priority column
# (10..0 -> 9..0, i.e. drops the old top row) and reshuffles the From/To/
# Drone text values, while also dropping the now-unused last row so the
# table shrinks from A1:D12 to A1:D11.
#
# Every text value needed after the edit already exists somewhere in the
# sheet before the edit - so instead of assigning new string literals
# (which this engine - like real Excel - would store as *numbers* when the
# text happens to look like an integer, or as a "quote-prefixed"/
# text-formatted cell when forced via NumberFormat, either way adding a
# style that is not part of the target), we relocate the existing text
# cells with Range.Copy. Copy preserves the shared-string cell type with
# no style changes, exactly like dragging/copy-pasting in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch row, well outside the table, used to "hold" a copy of each of
# the ten distinct text values (1..10) before any destination cells are
# overwritten - so the write order below can't clobber a value that is
# still needed as a source for a later cell.
$scratchRow = 100

# One source cell (row, column letter) per distinct text value, picked
# from the table as it exists before any edits.
$sources = @{
    "1"  = @(6,  "B")
    "2"  = @(4,  "B")
    "3"  = @(2,  "B")
    "4"  = @(2,  "C")
    "5"  = @(7,  "B")
    "6"  = @(8,  "B")
    "7"  = @(9,  "B")
    "8"  = @(11, "C")
    "9"  = @(6,  "C")
    "10" = @(4,  "C")
}

# Stage every distinct value into its own scratch column.
$scratchCol = @{}
$i = 1
foreach ($value in @("1","2","3","4","5","6","7","8","9","10")) {
    $rc = $sources[$value]
    $ws.Cells.Item($rc[0], $rc[1]).Copy($ws.Cells.Item($scratchRow, $i))
    $scratchCol[$value] = $i
    $i++
}

# Target table (row -> A, B/From, C/To, D/Drone) per the commit's fix.
$rows = @(
    @{ r = 2;  a = 9; b = "1"; c = "8"; d = "5" }
    @{ r = 3;  a = 8; b = "2"; c = "4"; d = "5" }
    @{ r = 4;  a = 7; b = "3"; c = "1"; d = "5" }
    @{ r = 5;  a = 6; b = "4"; c = "3"; d = "5" }
    @{ r = 6;  a = 5; b = "8"; c = "9"; d = "5" }
    @{ r = 7;  a = 4; b = "9"; c = "2"; d = "5" }
    @{ r = 8;  a = 3; b = "1"; c = "5"; d = "9" }
    @{ r = 9;  a = 2; b = "5"; c = "6"; d = "9" }
    @{ r = 10; a = 1; b = "6"; c = "7"; d = "9" }
    @{ r = 11; a = 0; b = "7"; c = "1"; d = "9" }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($scratchRow, $scratchCol[$row.b]).Copy($ws.Cells.Item($row.r, "B"))
    $ws.Cells.Item($scratchRow, $scratchCol[$row.c]).Copy($ws.Cells.Item($row.r, "C"))
    $ws.Cells.Item($scratchRow, $scratchCol[$row.d]).Copy($ws.Cells.Item($row.r, "D"))
}

# Wipe the scratch row now that every value has been relocated.
$ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, 10)).ClearContents()

# The old row 12 is no longer part of the (now 10-row) table - delete it
# so the sheet dimension shrinks from A1:D12 to A1:D11, matching the diff.
$ws.Rows(12).Delete()
